$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New "mixed"/"mixed2"/"intcol2" columns (K, L, M) and a new data row
# (row 6, "Quentin") -- mirrors the upstream fix that makes column type
# detection look at every cell in a column instead of just the first.
# Values are written in the same order the strings were first typed so
# that the shared-string table comes out in the expected order.
# ---------------------------------------------------------------------

# Column K ("mixed": number then text)
$ws.Range("K3").Value = "mixed"
$ws.Range("K4").Value = 123
$ws.Range("K5").Value = "abc"

# Column L ("mixed2": text then number)
$ws.Range("L3").Value = "mixed2"
$ws.Range("L4").Value = "abc"
$ws.Range("L5").Value = 123

# Row 6 ("Quentin") across the existing columns
$ws.Range("B6").Value = "Quentin"
$ws.Range("C6").Value = 125
$ws.Range("D6").Value = 12345679
$ws.Range("E6").Value = 12345678901
$ws.Range("F6").Value = 13.35
$ws.Range("G6").Formula = "=FALSE"
$ws.Range("H6").Value = 43913.913310185198
$ws.Range("I6").Formula = "=C6+F6"
$ws.Range("M6").Value = 1234

# Column M ("intcol2")
$ws.Range("M3").Value = "intcol2"
$ws.Range("M5").Value = 1234

# ---------------------------------------------------------------------
# Formatting: reuse the existing styles already present in the workbook
# (copy/paste-special "formats only" so no new style entries are added)
# ---------------------------------------------------------------------

# Header style (row 3) for the new K/L/M headers
$ws.Range("J3").Copy() | Out-Null
$ws.Range("K3:M3").PasteSpecial(-4122) | Out-Null

# Numeric style used by column A placeholder cells, for K4
$ws.Range("A4").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null

# Row 6 data cells take on the same per-column styling as row 5
$ws.Range("D5").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$ws.Range("F5").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
$ws.Range("G5").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null
$ws.Range("H5").Copy() | Out-Null
$ws.Range("H6").PasteSpecial(-4122) | Out-Null
$ws.Range("I5").Copy() | Out-Null
$ws.Range("I6").PasteSpecial(-4122) | Out-Null
$ws.Range("J3").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null

# Row 6 is very slightly shorter than the other data rows
$ws.Rows.Item(6).RowHeight = 15.5

# ---------------------------------------------------------------------
# Currency number format: "#,##0.00 €" -> "#,##0.00 $"
# ---------------------------------------------------------------------
$ws.Range("A5").ClearFormats() | Out-Null
$ws.Range("A5").NumberFormat = '#,##0.00\ "$"'

# ---------------------------------------------------------------------
# Selection left where the user ended up after entering the data
# ---------------------------------------------------------------------
$ws.Range("M6").Select() | Out-Null
